$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.561.55"
$ws.Range("E2").Value = "  -2.09%  "
$ws.Range("D3").Value = "3.026.31"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.11"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.56"
$ws.Range("E6").Value = "  -5.53%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "3.024.77"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("E10").Value = "  -4.66%  "
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("E12").Value = "  -2.97%  "
$ws.Range("E13").Value = "  -4.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.16"
$ws.Range("E14").Value = "  -6.09%  "
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("D16").Value = "3.529.19"
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("D17").Value = "62.555.09"
$ws.Range("E17").Value = "  -2.02%  "
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "3.026.13"
$ws.Range("E19").Value = "  -2.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.20"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.03"
$ws.Range("E21").Value = "  -3.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.689"
$ws.Range("E22").Value = "  -3.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.37"
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.02"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.45"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.31"
$ws.Range("E27").Value = "  +0.97%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  -4.23%  "
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.53"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -6.89%  "
$ws.Range("E35").Value = "  -1.27%  "
$ws.Range("D36").Value = "0.0₃0794"
$ws.Range("E36").Value = "  -7.70%  "
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.14"
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.25"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.02"
$ws.Range("E40").Value = "  -3.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  -17.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "420.69"
$ws.Range("E42").Value = "  -7.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.112"
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.279"
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").Value = "2.791.63"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.24"
$ws.Range("E47").Value = "  -8.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.66"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.33"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("E51").Value = "  -1.88%  "
